$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Week 4 tools cell: "Data Formats; readr; tidyr" -> add "surveys, "
# ------------------------------------------------------------------
$d.Content.Find.Execute("Data Formats; readr; tidyr", $true, $false, $false, $false, $false, $true, 1, $false, "Data Formats; surveys, readr; tidyr", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Week 11 tools cell: "tidycensus" -> "tidyverse mapping; color scales; projection"
#    (must happen before the week-12 tools edit, and before the topic
#    text "Census Data" is turned into a hyperlink)
# ------------------------------------------------------------------
$d.Content.Find.Execute("tidycensus", $true, $false, $false, $false, $false, $true, 1, $false, "tidyverse mapping; color scales; projection", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Week 12 tools cell: "sf; tigris; mapgl; mapbox; osm" -> "tidycensus; sf; tigris; mapgl; mapbox; osm"
# ------------------------------------------------------------------
$d.Content.Find.Execute("sf; tigris; mapgl; mapbox; osm", $true, $false, $false, $false, $false, $true, 1, $false, "tidycensus; sf; tigris; mapgl; mapbox; osm", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Turn three plain-text topic cells into hyperlinks (Hyperlink style)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Social Networks & Network Data") | Out-Null
$d.Hyperlinks.Add($r, "assignments/assign11.qmd") | Out-Null

$r = $d.Content
$r.Find.Execute("Census Data") | Out-Null
$d.Hyperlinks.Add($r, "assignments/assign12.qmd", "", "", "Cartography") | Out-Null

$r = $d.Content
$r.Find.Execute("Maps & GIS") | Out-Null
$d.Hyperlinks.Add($r, "assignments/assign13.qmd", "", "", "Geographic Data") | Out-Null

# ------------------------------------------------------------------
# 5) Second table ("Advanced Topics"): clear the week numbers and
#    insert a new "Missing Data" row after the "Text Data & Data
#    Scraping" row.
# ------------------------------------------------------------------
$t2 = $d.Tables(2)

# Clear week-number cells (row 2..5 hold 16,17,18,19 in column 1)
$weekNumbers = @("16", "17", "18", "19")
for ($i = 2; $i -le 5; $i++) {
    $cellRange = $t2.Rows($i).Cells(1).Range
    $cellRange.Find.Execute($weekNumbers[$i - 2], $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
}

# Insert the new row right after the "16 / Text Data & Data Scraping" row (row 2)
$newRow = $t2.Rows.Add($t2.Rows(3))
$newRow.Cells(2).Range.Text = "Missing Data"

Write-Host "done"
